$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "ArtistListPath" entry, recorded as the next row in the path list.
$ws.Range("A16").Value = "ArtistListPath"
$ws.Range("B16").Value = "Resources/DataLists/ArtistList.csv"

# Leave the selection where the author left it after adding the row.
$ws.Range("A18").Select() | Out-Null

# Window was un-minimized and moved/resized on save.
$excel.WindowState = -4143
$win = $wb.Windows.Item(1)
$win.Left = 7200
$win.Top = 1530
$win.Width = 21600
$win.Height = 11385
